$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row is inserted above the existing row 83,
# pushing rows 83..188 down to 84..189 (dimension grows from A1:T188 to
# A1:T189). Insert a whole row so everything below shifts down first.
$ws.Rows.Item(83).Insert()

# Populate the newly-inserted row 83 with the new record.
$ws.Range("A83").Value = 7
$ws.Range("B83").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C83").Value = "Ñuble"
$ws.Range("D83").Value = 44679
$ws.Range("E83").Value = 16
$ws.Range("F83").Value = "Fruta"
$ws.Range("G83").Value = 100104
$ws.Range("H83").Value = "Frutos de pepita"
$ws.Range("I83").Value = 100104005
$ws.Range("J83").Value = "Pera"
$ws.Range("K83").Value = "Packham's Triumph"
$ws.Range("L83").Value = "Primera"
$ws.Range("M83").Value = 160
$ws.Range("N83").Value = 8000
$ws.Range("O83").Value = 9000
$ws.Range("P83").Value = 8500
$ws.Range("Q83").Value = '$/caja 16 kilos empedrada'
$ws.Range("R83").Value = "Provincia de Curicó"
$ws.Range("S83").Value = 531
$ws.Range("T83").Value = 16

# Make sure the date cell keeps the date number format used throughout
# column D (style index 2 / numFmt 165 in the original workbook).
$ws.Range("D83").NumberFormat = "YYYY-MM-DD HH:MM:SS"
